$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1328
$ws1.Range("F3").Value = 1203
$ws1.Range("F4").Value = 900
$ws1.Range("F7").Value = 666
$ws1.Range("F8").Value = 110
$ws1.Range("F11").Value = 2403
$ws1.Range("F12").Value = 1598
$ws1.Range("F13").Value = 1421
$ws1.Range("F15").Value = 239
$ws1.Range("F16").Value = 576
$ws1.Range("F17").Value = 771
$ws1.Range("F18").Value = 58
$ws1.Range("F19").Value = 299
$ws1.Range("F24").Value = 4821
$ws1.Range("F25").Value = 214
$ws1.Range("F26").Value = 361
$ws1.Range("F27").Value = 65
$ws1.Range("F28").Value = 154
$ws1.Range("F31").Value = 95
$ws1.Range("F32").Value = 24
$ws1.Range("F34").Value = 702
$ws1.Range("F36").Value = 43
$ws1.Range("F37").Value = 238
$ws1.Range("F38").Value = 381
$ws1.Range("F39").Value = 1020
$ws1.Range("F41").Value = 101
$ws1.Range("F42").Value = 158
$ws1.Range("F43").Value = 124

# Sheet 4: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1328
$ws4.Range("F5").Value = 1203
$ws4.Range("F6").Value = 900
$ws4.Range("F11").Value = 667
$ws4.Range("F12").Value = 110
$ws4.Range("F17").Value = 2403
$ws4.Range("F18").Value = 1598
$ws4.Range("F19").Value = 1421
$ws4.Range("F21").Value = 239
$ws4.Range("F22").Value = 576
$ws4.Range("F24").Value = 771
$ws4.Range("F25").Value = 58
$ws4.Range("F26").Value = 299
$ws4.Range("F29").Value = 4821
$ws4.Range("F30").Value = 214
$ws4.Range("F31").Value = 361
$ws4.Range("F32").Value = 65
$ws4.Range("F33").Value = 154
$ws4.Range("F36").Value = 95
$ws4.Range("F37").Value = 24
$ws4.Range("F39").Value = 702
$ws4.Range("F40").Value = 43
$ws4.Range("F41").Value = 381
$ws4.Range("F42").Value = 1020
$ws4.Range("F44").Value = 101
$ws4.Range("F45").Value = 158
$ws4.Range("F46").Value = 124
